# Insert a new weekly price record at row 61 for "Poroto verde" (Feria
# Lagunitas de Puerto Montt). Inserting the row pushes the existing
# rows 61-75 down to 62-76, preserving all of their data/formatting,
# matching the source diff exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(61).Insert()

$newDate = Get-Date -Year 2022 -Month 5 -Day 13 -Hour 0 -Minute 0 -Second 0

$ws.Cells.Item(61, 1).Value = 4
$ws.Cells.Item(61, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(61, 3).Value = "Los Lagos"
$ws.Cells.Item(61, 4).Value = $newDate
$ws.Cells.Item(61, 5).Value = 10
$ws.Cells.Item(61, 6).Value = 100112031
$ws.Cells.Item(61, 7).Value = "Poroto verde"
$ws.Cells.Item(61, 8).Value = "Magnum"
$ws.Cells.Item(61, 9).Value = "Primera"
$ws.Cells.Item(61, 10).Value = 45
$ws.Cells.Item(61, 11).Value = 30000
$ws.Cells.Item(61, 12).Value = 30000
$ws.Cells.Item(61, 13).Value = 30000
$ws.Cells.Item(61, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(61, 15).Value = "Región Metropolitana"
$ws.Cells.Item(61, 16).Value = 1200
$ws.Cells.Item(61, 17).Value = 25
$ws.Cells.Item(61, 18).Value = "Hortaliza"
